# Deploying to gh-pages from  @ 0d0ffeaff9e0087f25ede5ade66aaa8603c52931 🚀
# Adds a new "2023" data column (column Q) to the sheet, mirroring the
# formatting of the existing 2022 column (P) and refreshing row heights
# / column widths the way Excel does when a sheet is touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new year column Q, copying each row's formatting from
#        the corresponding cell in column P (the previous last year) and
#        then writing in the new value. ---

# Row 3 : header year
$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("Q3").Value = 2023

# Row 5 : total (Кыргызская Республика)
$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 2385.9

# Row 6 : Баткенская область
$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 112.1

# Row 7 : Джалал-Абадская область
$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("Q7").Value = 267.89999999999998

# Row 8 : Иссык-Кульская область
$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("Q8").Value = 230.9

# Row 9 : Нарынская область
$ws.Range("P9").Copy($ws.Range("Q9"))
$ws.Range("Q9").Value = 249.7

# Row 10 : Ошская область
$ws.Range("P10").Copy($ws.Range("Q10"))
$ws.Range("Q10").Value = 287

# Row 11 : Таласская область
$ws.Range("P11").Copy($ws.Range("Q11"))
$ws.Range("Q11").Value = 334.7

# Row 12 : Чуйская область
$ws.Range("P12").Copy($ws.Range("Q12"))
$ws.Range("Q12").Value = 851

# Row 13 : г.Бишкек
$ws.Range("P13").Copy($ws.Range("Q13"))
$ws.Range("Q13").Value = 48.5

# Row 14 : г.Ош
$ws.Range("P14").Copy($ws.Range("Q14"))
$ws.Range("Q14").Value = 4.2

# Row 15 : section header (percentage table) - formatting only, no value
$ws.Range("P15").Copy($ws.Range("Q15"))

# Row 16 : Кыргызская Республика (%)
$ws.Range("P16").Copy($ws.Range("Q16"))
$ws.Range("Q16").Value = 26.890545708088244

# Row 17 : Баткенская область (%)
$ws.Range("P17").Copy($ws.Range("Q17"))
$ws.Range("Q17").Value = 15.490056759274875

# Row 18 : Джалал-Абадская область (%)
$ws.Range("P18").Copy($ws.Range("Q18"))
$ws.Range("Q18").Value = 22.218388220841799

# Row 19 : Иссык-Кульская область (%)
$ws.Range("P19").Copy($ws.Range("Q19"))
$ws.Range("Q19").Value = 29.614327895683314

# Row 20 : Нарынская область (%)
$ws.Range("P20").Copy($ws.Range("Q20"))
$ws.Range("Q20").Value = 30.104452089276922

# Row 21 : Ошская область (%)
$ws.Range("P21").Copy($ws.Range("Q21"))
$ws.Range("Q21").Value = 21.825966598728439

# Row 22 : Таласская область (%)
$ws.Range("P22").Copy($ws.Range("Q22"))
$ws.Range("Q22").Value = 32.351574864874735

# Row 23 : Чуйская область (%)
$ws.Range("P23").Copy($ws.Range("Q23"))
$ws.Range("Q23").Value = 30.810022297218843

# Row 24 : г.Бишкек (%)
$ws.Range("P24").Copy($ws.Range("Q24"))
$ws.Range("Q24").Value = 29.193884213235311

# Row 25 : г.Ош (%)
$ws.Range("P25").Copy($ws.Range("Q25"))
$ws.Range("Q25").Value = 7.4362892319581295

# --- 2. Refresh the explicit row heights for rows 4-25 (Excel marks
#        these as custom-height once the row is touched by the new
#        column, even though the value matches the sheet default). ---
$ws.Rows("4:25").RowHeight = 15

# --- 3. Narrow columns A:C slightly, matching the width Excel computed
#        for the table after the new column was added. ---
$ws.Columns("A:C").ColumnWidth = 35.45

# --- 4. Reset the active selection back to the top-left cell (the
#        workbook was re-saved with the default selection rather than
#        the ad-hoc S3 selection left over from editing). ---
$ws.Range("A1").Select()
